# edit.ps1 - PowerPoint COM-interop script reproducing the commit's changes:
#   1. Update the "datetimeFigureOut" date placeholder text from 4/1/2020 to
#      4/5/2020 on the slide master and every slide layout.
#   2. On slide 14, expand the "content placeholder" text with additional
#      Hebrew sentences, and turn on "shrink text on overflow" (normAutofit).
#   3. Rename shape id 4 ("מלבן 3" -> "כותרת 3") on slides 15, 16, 17 and 18.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder: "4/1/2020" -> "4/5/2020" on master + all layouts
# ---------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "4/1/2020") {
                $tr.Text = "4/5/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateShape $layouts.Item($L).Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 14: extend the "הגענו למסקנה " paragraph with the new
#    explanatory sentences, and enable shrink-text-on-overflow autofit.
# ---------------------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$contentShape = $slide14.Shapes.Item(2)

$contentShape.TextFrame.AutoSize = 2

$tr = $contentShape.TextFrame.TextRange
$targetPara = $tr.Paragraphs(2, 1)

# The paragraph currently ends with a single lone-space run ("הגענו למסקנה" + " ").
# Rewrite that trailing space run's text in place (keeps its own <a:rPr/>),
# then append the remaining new sentences as additional runs after it.
# NOTE: Paragraphs(...).Text includes the trailing paragraph-mark (CR), so the
# last *real* character is at index Length-1, not Length.
$paraLen = $targetPara.Text.Length
$lastChar = $targetPara.Characters($paraLen - 1, 1)
$lastChar.Text = " כי גם שימוש "

$run = $targetPara.InsertAfter("באוגמנטציות")
$run = $run.InsertAfter(" וגם הגדלת ")
$run = $run.InsertAfter("הדאטא")
$run = $run.InsertAfter(" סט לא השפיעו במידה רבה על דיוק המסווג. ההערכה שלנו היא שהגדלת מספר ")
$run = $run.InsertAfter("הטרנספורמתיות")
$run = $run.InsertAfter(" תשפר ")
$run = $run.InsertAfter("את דיוק המסווג.")

# ---------------------------------------------------------------------
# 3) Rename shape id 4 on slides 15, 16, 17, 18: "מלבן 3" -> "כותרת 3"
# ---------------------------------------------------------------------
foreach ($idx in 15, 16, 17, 18) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq 4 -and $sh.Name -eq "מלבן 3") {
            $sh.Name = "כותרת 3"
        }
    }
}
